# "modified the login code"
# Appends 5 new "seyntt logged in" event rows to the Sheet2 login log,
# continuing directly after the existing last row (row 123).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Find the last used row in the log (mirrors how the event-driven logger
# locates the next free row) instead of hard-coding 123, so the script is
# robust if run against a sheet whose data extends further.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# New login events to append: user name, message, login date (date-only
# serial) and the precise login timestamp (date + time serial).
$events = @(
  @{ User = "seyntt"; Message = "seyntt logged in"; LoginDate = 45787; LoginTime = 45787.746875 },
  @{ User = "seyntt"; Message = "seyntt logged in"; LoginDate = 45787; LoginTime = 45787.748912037 },
  @{ User = "seyntt"; Message = "seyntt logged in"; LoginDate = 45788; LoginTime = 45788.8391782407 },
  @{ User = "seyntt"; Message = "seyntt logged in"; LoginDate = 45788; LoginTime = 45788.8397222222 },
  @{ User = "seyntt"; Message = "seyntt logged in"; LoginDate = 45788; LoginTime = 45788.8403356481 }
)

# The most recently-used style for the Date/Time columns (picked up from
# the last existing row) so the new rows continue the same formatting.
$dateFormat = $ws.Cells.Item($lastRow, 3).NumberFormat
$timeFormat = $ws.Cells.Item($lastRow, 4).NumberFormat

$row = $lastRow
foreach ($evt in $events) {
  $row = $row + 1

  $ws.Cells.Item($row, 1).Value2 = $evt.User
  $ws.Cells.Item($row, 2).Value2 = $evt.Message

  $ws.Cells.Item($row, 3).Value2 = $evt.LoginDate
  $ws.Cells.Item($row, 3).NumberFormat = $dateFormat

  $ws.Cells.Item($row, 4).Value2 = $evt.LoginTime
  $ws.Cells.Item($row, 4).NumberFormat = $timeFormat

  # Match the row height ("ht") that every other row in this sheet carries.
  $ws.Rows.Item($row).RowHeight = 15
}

Write-Output "Appended $($events.Count) login rows ($($lastRow + 1):$row)"
